$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Formula = "'" + $val
    $cell.ClearFormats()
}

Set-TextCell 2 4 "20.529.83"
Set-TextCell 2 5 "  +1.70%  "

Set-TextCell 3 4 "1.469.68"
Set-TextCell 3 5 "  +2.16%  "

Set-TextCell 4 4 "1.007"
Set-TextCell 4 5 "  +0.17%  "

Set-TextCell 5 4 "0.9581"
Set-TextCell 5 5 "  +4.66%  "

Set-TextCell 6 4 "277.19"
Set-TextCell 6 5 "  +0.46%  "

Set-TextCell 7 4 "0.3596"
Set-TextCell 7 5 "  -0.80%  "

Set-TextCell 8 4 "0.3078"
Set-TextCell 8 5 "  +0.16%  "

Set-TextCell 9 2 "Polygon"
Set-TextCell 9 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 9 4 "1.079"
Set-TextCell 9 5 "  +5.37%  "

Set-TextCell 10 2 "OKB"
Set-TextCell 10 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 10 4 "39.50"
Set-TextCell 10 5 "  +1.66%  "

Set-TextCell 11 4 "0.06622"
Set-TextCell 11 5 "  +2.01%  "

Set-TextCell 12 4 "1.001"
Set-TextCell 12 5 "  +0.19%  "

Set-TextCell 13 4 "5.468"
Set-TextCell 13 5 "  +2.26%  "

Set-TextCell 14 4 "18.07"
Set-TextCell 14 5 "  +3.00%  "

Set-TextCell 15 4 "6.143"
Set-TextCell 15 5 "  +1.70%  "

Set-TextCell 16 4 "0.9585"
Set-TextCell 16 5 "  +2.09%  "

Set-TextCell 17 4 "0.00001020"
Set-TextCell 17 5 "  +1.03%  "

Set-TextCell 18 4 "1.466.27"
Set-TextCell 18 5 "  +2.38%  "

Set-TextCell 19 4 "0.05941"
Set-TextCell 19 5 "  +5.64%  "

Set-TextCell 20 4 "68.56"
Set-TextCell 20 5 "  +1.26%  "

Set-TextCell 21 5 "  +1.33%  "

Set-TextCell 22 4 "14.50"
Set-TextCell 22 5 "  +1.34%  "

Set-TextCell 23 4 "11.18"
Set-TextCell 23 5 "  +2.84%  "

Set-TextCell 24 4 "2.274"
Set-TextCell 24 5 "  +1.35%  "

Set-TextCell 25 4 "20.530.60"
Set-TextCell 25 5 "  +1.66%  "

Set-TextCell 26 4 "143.50"
Set-TextCell 26 5 "  +3.83%  "

Set-TextCell 27 4 "2.094"
Set-TextCell 27 5 "  -1.86%  "

Set-TextCell 28 4 "17.13"
Set-TextCell 28 5 "  +1.59%  "

Set-TextCell 29 4 "1.626.13"
Set-TextCell 29 5 "  +2.70%  "

Set-TextCell 30 4 "113.64"
Set-TextCell 30 5 "  +3.38%  "

Set-TextCell 31 4 "3.873"
Set-TextCell 31 5 "  +1.32%  "

Set-TextCell 32 2 "Filecoin"
Set-TextCell 32 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 32 4 "4.949"
Set-TextCell 32 5 "  +2.36%  "

Set-TextCell 33 2 "Stellar"
Set-TextCell 33 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 33 4 "0.07977"
Set-TextCell 33 5 "  +4.43%  "

Set-TextCell 34 4 "0.8009"
Set-TextCell 34 5 "  -0.99%  "

Set-TextCell 35 4 "1.223"
Set-TextCell 35 5 "  +8.41%  "

Set-TextCell 36 4 "1.455"
Set-TextCell 36 5 "  -0.68%  "

Set-TextCell 37 4 "0.05761"
Set-TextCell 37 5 "  -1.21%  "

Set-TextCell 38 4 "4.702"
Set-TextCell 38 5 "  +1.53%  "

Set-TextCell 39 4 "0.02042"
Set-TextCell 39 5 "  +2.75%  "

Set-TextCell 40 4 "0.9587"
Set-TextCell 40 5 "  +3.84%  "

Set-TextCell 41 4 "10.35"
Set-TextCell 41 5 "  +1.85%  "

Set-TextCell 42 4 "0.1861"
Set-TextCell 42 5 "  +1.29%  "

Set-TextCell 43 4 "7.299"
Set-TextCell 43 5 "  +1.58%  "

Set-TextCell 44 4 "0.5259"
Set-TextCell 44 5 "  +1.01%  "

Set-TextCell 45 4 "3.513"
Set-TextCell 45 5 "  +0.55%  "

Set-TextCell 46 4 "12.10"
Set-TextCell 46 5 "  +2.46%  "

Set-TextCell 47 4 "118.81"
Set-TextCell 47 5 "  +1.84%  "

Set-TextCell 48 4 "0.5181"
Set-TextCell 48 5 "  +1.98%  "

Set-TextCell 49 4 "1.802"
Set-TextCell 49 5 "  +3.41%  "

Set-TextCell 50 4 "0.06441"

Set-TextCell 51 4 "0.9895"
Set-TextCell 51 5 "  +0.10%  "

